$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update URL row (row 2) ---
# Old URL moves to new column D, new URL goes into B and C (both hyperlinked)
$ws.Range("D2").Value2 = "https://hecmoptiam.com/Home"
$ws.Range("B2").Value2 = "https://smartnsc.com/"
$ws.Range("C2").Value2 = "https://smartnsc.com/"

# --- UserName row (row 3) stays the same, copied across C and D ---
$ws.Range("C3").Value2 = "Amitthakur"
$ws.Range("D3").Value2 = "Amitthakur"

# --- Password row (row 4) ---
# Old password moves to new column D, new password goes into B and C (both hyperlinked)
$ws.Range("D4").Value2 = "Aamit5555500000@@"
$ws.Range("B4").Value2 = "Aamit55555000@"
$ws.Range("C4").Value2 = "Aamit55555000@"

# --- CaseNumber rows (5 and 6) get new case numbers ---
$ws.Range("B5").Value2 = "7744000279"
$ws.Range("B6").Value2 = "7755000008"

# --- Hyperlinks: add in the same order so relationship ids line up ---
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:Aamit55555000@") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://smartnsc.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://smartnsc.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:Aamit55555000@") | Out-Null

# Re-apply the workbook's existing "Hyperlink" cell style so hyperlinked /
# hyperlink-styled cells reuse the original style index instead of a new one
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("C4").Style = "Hyperlink"

# --- New column widths for C and D (closest achievable to authored 21.140625 / 30.140625) ---
$ws.Columns("C").ColumnWidth = 20.3
$ws.Columns("D").ColumnWidth = 29.3

# --- Selection matches the authored workbook ---
$ws.Range("D13").Select() | Out-Null
